$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (shifts NIF, Colegio, Mesa one column right)
$ws.Columns.Item(2).Insert()

# Remove the old "Mesa" column, now shifted to column E
$ws.Columns.Item(5).Delete()

# Fill the new column B with the Email header and values
$ws.Range("B1").Value = "Email"
$ws.Range("B4").Value = "email3"
$ws.Range("B2").Value = "email1"
$ws.Range("B3").Value = "email2"

# Move the active selection to B5
$ws.Range("B5").Select()
